$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price column (D) ---
# Values are textual (may look numeric); force text format so
# Excel does not coerce them into numbers and lose formatting,
# then restore the default "Normal" style so no residual number
# format style is left behind on the cell.
$priceUpdates = @{
    "D2" = "61.088.20"
    "D3" = "2.925.54"
    "D5" = "593.01"
    "D6" = "145.07"
    "D9" = "6.98"
    "D11" = "0.440"
    "D13" = "33.73"
    "D15" = "3.410.85"
    "D16" = "61.066.46"
    "D17" = "6.74"
    "D18" = "2.927.91"
    "D19" = "434.99"
    "D22" = "7.10"
    "D23" = "81.44"
    "D25" = "2.20"
    "D26" = "11.86"
    "D29" = "2.60"
    "D30" = "6.96"
    "D32" = "26.67"
    "D34" = "0.0₃0868"
    "D40" = "8.59"
    "D41" = "42.01"
    "D42" = "0.285"
    "D43" = "375.05"
    "D44" = "0.0347"
    "D45" = "2.683.79"
    "D46" = "134.12"
    "D48" = "23.93"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# --- Update Volume(1h) column (E) ---
$volumeUpdates = @{
    "E2" = "  +0.53%  "
    "E3" = "  +0.60%  "
    "E4" = "  +0.02%  "
    "E5" = "  +0.61%  "
    "E6" = "  +0.09%  "
    "E7" = "  +0.03%  "
    "E8" = "  +0.18%  "
    "E9" = "  +4.46%  "
    "E10" = "  -0.78%  "
    "E11" = "  -0.94%  "
    "E12" = "  -0.40%  "
    "E13" = "  +0.51%  "
    "E14" = "  -0.09%  "
    "E15" = "  +0.69%  "
    "E16" = "  +0.53%  "
    "E17" = "  +0.48%  "
    "E18" = "  +0.74%  "
    "E19" = "  +1.20%  "
    "E20" = "  -0.52%  "
    "E21" = "  -0.72%  "
    "E22" = "  +0.29%  "
    "E23" = "  +0.01%  "
    "E24" = "  +1.93%  "
    "E25" = "  -0.94%  "
    "E26" = "  -0.56%  "
    "E27" = "  -0.04%  "
    "E28" = "  +0.62%  "
    "E29" = "  -0.49%  "
    "E30" = "  -1.89%  "
    "E31" = "  +3.36%  "
    "E32" = "  +0.70%  "
    "E33" = "  +0.10%  "
    "E34" = "  +1.94%  "
    "E35" = "  -0.13%  "
    "E36" = "  +0.74%  "
    "E37" = "  -1.16%  "
    "E38" = "  -0.31%  "
    "E39" = "  -0.30%  "
    "E40" = "  -0.24%  "
    "E41" = "  +2.71%  "
    "E42" = "  -3.42%  "
    "E43" = "  +0.37%  "
    "E44" = "  -1.08%  "
    "E45" = "  -0.52%  "
    "E46" = "  +1.42%  "
    "E48" = "  -0.57%  "
    "E49" = "  -0.83%  "
    "E50" = "  -1.94%  "
    "E51" = "  -0.16%  "
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
